$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
try {
  $cs = $nm.ColorScheme
  Write-Output ("count=" + $cs.Count)
} catch {
  Write-Output ("ERR: " + $_.Exception.Message)
}
